$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2020" column (Q) to the table -----------------------------
# Header cell: year 2020, formatted like the existing year header in P3.
$ws.Range("Q3").Value = 2020
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)   # xlPasteFormats

# Data cell: the 2020 value for the agriculture orientation index,
# formatted like the existing data cell in P4.
$ws.Range("Q4").Value = 0.067156049127444606
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)   # xlPasteFormats

# --- Re-format the data row's numbers --------------------------------
# The row used a custom "0.0" format; switch the whole row (now including
# the new Q4 cell) to the standard "0.00" number format.
$ws.Range("D4:Q4").NumberFormat = "0.00"

# --- Reset the sheet selection back to the top-left cell -----------------
$ws.Range("A1").Select()
